# Generate Report for Archive
# 1) Status text changes from "Ready for handoff" to "In Translation"
#    (Overview!E2/F2, zh-cn!C2, de-de!C2 — all shared the same "Ready for
#    handoff" string).
# 2) Narrow the "status" columns (Overview E & F, zh-cn C, de-de C) from
#    ~17.22 down to ~13.41 characters to fit the shorter "In Translation"
#    text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh       = $wb.Worksheets.Item("zh-cn")
$wsDe       = $wb.Worksheets.Item("de-de")

# --- Update the status values -------------------------------------------
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZh.Range("C2").Value = "In Translation"
$wsDe.Range("C2").Value = "In Translation"

# --- Narrow the status columns -------------------------------------------
$wsOverview.Range("E1:F1").EntireColumn.ColumnWidth = 12.5
$wsZh.Range("C1").EntireColumn.ColumnWidth = 12.5
$wsDe.Range("C1").EntireColumn.ColumnWidth = 12.5
